$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (rows 2-25), columns B,C,D,E,F,I,J,K,L,M,N
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.056668791459715; "D" = 1.06877599645973; "E" = 1.052691140899335; "F" = 1.075244493472078; "I" = 1.054580223426608; "J" = 1.061669224238124; "K" = 1.071480274395791; "L" = 1.05543917352213; "M" = 1.077931574523797; "N" = 1.063176916311704 }
    3 = @{ "B" = 1.02; "C" = 1.058187082097106; "D" = 1.069837089383537; "E" = 1.054026462375493; "F" = 1.07655252367017; "I" = 1.055104899635551; "J" = 1.062836940710988; "K" = 1.072356801852512; "L" = 1.056586051399844; "M" = 1.079055680242195; "N" = 1.064346291075899 }
    4 = @{ "B" = 1.02; "C" = 1.059167296141315; "D" = 1.07052240456511; "E" = 1.054888035637026; "F" = 1.077397783988321; "I" = 1.055442185329456; "J" = 1.063589908028761; "K" = 1.072922090208696; "L" = 1.057325154484324; "M" = 1.0797813769249; "N" = 1.065100327693629 }
    5 = @{ "B" = 1.02; "C" = 1.059578855919119; "D" = 1.070810208503757; "E" = 1.055249658342478; "F" = 1.077752868634738; "I" = 1.055583453371064; "J" = 1.063905834959876; "K" = 1.07315929094637; "L" = 1.057635162144732; "M" = 1.080086064436976; "N" = 1.065416703277224 }
    6 = @{ "B" = 1.02; "C" = 1.05964792827298; "D" = 1.070858514440859; "E" = 1.0553103424594; "F" = 1.077812473672088; "I" = 1.055607142099205; "J" = 1.06395884431055; "K" = 1.073199091939814; "L" = 1.057687172310785; "M" = 1.080137199787901; "N" = 1.065469787907256 }
    7 = @{ "B" = 1.02; "C" = 1.059172797463226; "D" = 1.070526251398669; "E" = 1.054892869936529; "F" = 1.077402529671949; "I" = 1.055444075025991; "J" = 1.063594131888882; "K" = 1.072925261445494; "L" = 1.057329299607069; "M" = 1.079785449720992; "N" = 1.065104557552117 }
    8 = @{ "B" = 1.02; "C" = 1.05718237068587; "D" = 1.069134866171341; "E" = 1.053142935341904; "F" = 1.075686784475096; "I" = 1.054758000863465; "J" = 1.062064407193245; "K" = 1.071776893626905; "L" = 1.055827392888272; "M" = 1.078311821491944; "N" = 1.063572660471949 }
    9 = @{ "B" = 1.02; "C" = 1.053657536915863; "D" = 1.066673036215274; "E" = 1.050040039428788; "F" = 1.072654554191685; "I" = 1.053531934376027; "J" = 1.059348407711782; "K" = 1.069738697405493; "L" = 1.053157486885594; "M" = 1.075702020366441; "N" = 1.060852803959693 }
    10 = @{ "B" = 1.02; "C" = 1.051295308726367; "D" = 1.065024780622808; "E" = 1.047957943513224; "F" = 1.070626724937608; "I" = 1.052702833544543; "J" = 1.05752353148338; "K" = 1.068369784531649; "L" = 1.051361354850346; "M" = 1.073953003139323; "N" = 1.059025336197731 }
    11 = @{ "B" = 1.02; "C" = 1.050269375117456; "D" = 1.064309333321924; "E" = 1.047053052465871; "F" = 1.069747062681743; "I" = 1.052340994729449; "J" = 1.056729863216388; "K" = 1.067774567356987; "L" = 1.050579658418299; "M" = 1.073193409181385; "N" = 1.058230540830765 }
    12 = @{ "B" = 1.02; "C" = 1.049887824170208; "D" = 1.064043316903167; "E" = 1.046716425049979; "F" = 1.069420069289574; "I" = 1.052206161859452; "J" = 1.056434526677642; "K" = 1.067553100708999; "L" = 1.050288697109295; "M" = 1.072910915945066; "N" = 1.057934784880258 }
    13 = @{ "B" = 1.02; "C" = 1.049969689762352; "D" = 1.064100390523126; "E" = 1.046788656045539; "F" = 1.069490221810056; "I" = 1.052235103486756; "J" = 1.056497901627803; "K" = 1.067600623183618; "L" = 1.050351136805297; "M" = 1.072971527483271; "N" = 1.057998249830118 }
    14 = @{ "B" = 1.02; "C" = 1.050237845725593; "D" = 1.064287349810897; "E" = 1.047025237214639; "F" = 1.069720038401615; "I" = 1.052329858195316; "J" = 1.056705461528874; "K" = 1.067756268567413; "L" = 1.050555619852106; "M" = 1.073170065305779; "N" = 1.058206104490055 }
    15 = @{ "B" = 1.02; "C" = 1.050403002363022; "D" = 1.064402505987406; "E" = 1.047170934712272; "F" = 1.069861602946341; "I" = 1.052388182615136; "J" = 1.05683327515403; "K" = 1.067852116769566; "L" = 1.050681528222224; "M" = 1.07329234492318; "N" = 1.05833409962522 }
    16 = @{ "B" = 1.02; "C" = 1.051363330806479; "D" = 1.065072225317138; "E" = 1.048017927117743; "F" = 1.070685070904256; "I" = 1.052726787552267; "J" = 1.057576130423515; "K" = 1.068409234645074; "L" = 1.051413149231064; "M" = 1.074003366721282; "N" = 1.059078009834393 }
    17 = @{ "B" = 1.02; "C" = 1.051964888376488; "D" = 1.065491852196431; "E" = 1.048548324337379; "F" = 1.071201177307862; "I" = 1.052938424109799; "J" = 1.058041164418215; "K" = 1.068758035070301; "L" = 1.051871008969096; "M" = 1.074448762479639; "N" = 1.059543704230707 }
    18 = @{ "B" = 1.02; "C" = 1.0523154710123; "D" = 1.065736445794432; "E" = 1.048857375489693; "F" = 1.071502059687914; "I" = 1.053061595109414; "J" = 1.058312075274828; "K" = 1.068961246626787; "L" = 1.052137689435101; "M" = 1.07470833690595; "N" = 1.05981499981181 }
    19 = @{ "B" = 1.02; "C" = 1.052434960827993; "D" = 1.065819817625231; "E" = 1.048962699882341; "F" = 1.07160462694911; "I" = 1.053103547055343; "J" = 1.058404392230032; "K" = 1.069030496343321; "L" = 1.052228556183693; "M" = 1.074796808415312; "N" = 1.059907447867679 }
    20 = @{ "B" = 1.02; "C" = 1.051900377615161; "D" = 1.065446847587282; "E" = 1.048491450968374; "F" = 1.071145819932602; "I" = 1.052915745780941; "J" = 1.057991305451734; "K" = 1.068720636725915; "L" = 1.051821924468931; "M" = 1.074400998265596; "N" = 1.059493774458774 }
    21 = @{ "B" = 1.02; "C" = 1.050158893645014; "D" = 1.064232302389564; "E" = 1.04695558413519; "F" = 1.06965237004043; "I" = 1.052301967194859; "J" = 1.056644355094089; "K" = 1.067710445319132; "L" = 1.050495421454648; "M" = 1.073111610479278; "N" = 1.058144911277125 }
    22 = @{ "B" = 1.02; "C" = 1.049061207956447; "D" = 1.063467119421177; "E" = 1.045986964395054; "F" = 1.068711942166021; "I" = 1.051913571280373; "J" = 1.055794386646723; "K" = 1.067073116849247; "L" = 1.04965789439103; "M" = 1.072298914573113; "N" = 1.057293735777071 }
    23 = @{ "B" = 1.02; "C" = 1.049643376129396; "D" = 1.063872906268571; "E" = 1.04650073202613; "F" = 1.069210619462926; "I" = 1.052119704655062; "J" = 1.056245266753048; "K" = 1.067411185340063; "L" = 1.050102218504634; "M" = 1.072729932604833; "N" = 1.057745256184863 }
    24 = @{ "B" = 1.02; "C" = 1.05192952816519; "D" = 1.065467183751244; "E" = 1.048517150586379; "F" = 1.071170834022972; "I" = 1.052925993985276; "J" = 1.058013835609611; "K" = 1.068737536162992; "L" = 1.051844104820216; "M" = 1.074422581529825; "N" = 1.05951633661206 }
    25 = @{ "B" = 1.02; "C" = 1.054570920863576; "D" = 1.067310696822605; "E" = 1.050844552440964; "F" = 1.073439548214332; "I" = 1.053850952337446; "J" = 1.060053028233541; "K" = 1.070267382039169; "L" = 1.053850539653869; "M" = 1.076378302954009; "N" = 1.061558425123414 }
}

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}

Write-Host "Updated vm_pu values for rows 2-25"